$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B-column values for rows 2..127 (row 2 is index 0 of this array)
$bValues = @(856,516,557,579,590,767,805,806,807,829,513,517,534,535,544,547,551,569,582,591,592,600,611,646,647,649,650,653,654,655,661,662,677,705,706,709,715,719,720,727,730,733,734,735,737,742,743,744,748,753,756,758,760,763,768,770,771,772,774,775,787,797,798,813,821,822,823,831,832,839,843,844,848,849,853,510,574,576,581,583,651,671,675,686,687,713,716,718,728,739,747,749,750,751,752,754,762,769,773,777,784,790,791,793,796,801,803,809,810,834,512,603,683,700,702,717,738,764,779,782,802,740,778,828,833,660)

for ($i = 0; $i -lt $bValues.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

# Add the new rows 113..127 in column A with the corresponding serial labels (U7-4 .. U7-18)
for ($row = 113; $row -le 127; $row++) {
    $n = $row - 113 + 4
    $ws.Cells.Item($row, 1).Value = "U7-" + $n
}
